$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 15 and 16 were still blank template rows; fill them in with the new
# --- "no conformidad" entries for the week of 8-14 April 2016, matching the
# --- look & feel of the rows above (copy formatting from row 14 first).
$ws.Range("A14:G14").Copy()
$ws.Range("A15:G15").PasteSpecial(-4122)
$ws.Range("A14:G14").Copy()
$ws.Range("A16:G16").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("B15").Value = "Las tareas atrasadas no cuentan con un comentario de seguimiento."
$ws.Range("C15").Value = "Ventas"
$ws.Range("D15").Value = (Get-Date -Year 2016 -Month 4 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E15").Value = (Get-Date -Year 2016 -Month 4 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F15").Value = "Cerrada"
$ws.Range("G15").Value = "Esta anomalia solo se reporta ya que son tareas ya realizadas. Se presentan para futuras actividades y se escalan."
$ws.Rows("15").RowHeight = 75

$ws.Range("B16").Value = "La actividade del 14 no tiene comentario."
$ws.Range("C16").Value = "Compras"
$ws.Range("D16").Value = (Get-Date -Year 2016 -Month 4 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E16").Value = (Get-Date -Year 2016 -Month 4 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F16").Value = "Cerrada"
$ws.Range("G16").Value = "Esta anomalia solo se reporta ya que son tareas ya realizadas. Se presentan para futuras actividades y se escalan."
$ws.Rows("16").RowHeight = 75

# --- Move the cursor / viewport to reflect where the audit work left off.
[void]$ws.Range("C15").Select()
$excel.ActiveWindow.ScrollRow = 13
